# Bank statement content update: new cardholder, new card number, new
# transaction period (Oct 2024 instead of Oct 2023) with updated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: cardholder name
$ws.Range("C2").Value = "Hartmut"

# Card number (must stay text, not become a number) + surname
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)   # xlPasteValues - collapses formula to a plain string value
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 11.10.2024"

# Transaction row 6
$ws.Range("B6").Value = "12.10."
$ws.Range("C6").Value = "13.10."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "75,63-"

# Transaction row 7
$ws.Range("B7").Value = "16.10."
$ws.Range("C7").Value = "17.10."
$ws.Range("D7").Value = "MCDONALDS Tecklenburg"
$ws.Range("E7").Value = "8,72-"

# Transaction row 8
$ws.Range("B8").Value = "20.10."
$ws.Range("C8").Value = "21.10."
$ws.Range("D8").Value = "KARTENZ./20.10 LIDL RO"
$ws.Range("E8").Value = "17,16-"

# Transaction row 9
$ws.Range("B9").Value = "23.10."
$ws.Range("C9").Value = "24.10."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,89-"

# Row 10 no longer holds a transaction - clear it out entirely.
$ws.Range("B10").Value = $null
$ws.Range("C10").Value = $null
$ws.Range("D10").Value = $null
$ws.Range("E10").Value = $null
# ... and match its blank-row formatting (style 12) to the row above it.
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)   # xlPasteFormats

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 28.10.2024"
$ws.Range("E12").Value = "126,40-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.11.2024"
